# Update the division problems in the practice table.
# Each data row (1, 5, 9, 13, 17) of Table 1 holds 5 problems (columns 1-5).
# Replacements are applied per-cell by position (not via document-wide
# Find/Replace) because a couple of the original values ("55÷8=") repeat
# with different replacements, so a global text replace would be ambiguous.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$edits = @(
    @{ Row = 1;  Col = 1; New = "43÷8=" },
    @{ Row = 1;  Col = 2; New = "66÷7=" },
    @{ Row = 1;  Col = 3; New = "96÷4=" },
    @{ Row = 1;  Col = 4; New = "69÷4=" },
    @{ Row = 1;  Col = 5; New = "23÷3=" },

    @{ Row = 5;  Col = 1; New = "13÷5=" },
    @{ Row = 5;  Col = 2; New = "10÷7=" },
    @{ Row = 5;  Col = 3; New = "26÷4=" },
    @{ Row = 5;  Col = 4; New = "32÷4=" },
    @{ Row = 5;  Col = 5; New = "84÷9=" },

    @{ Row = 9;  Col = 1; New = "56÷9=" },
    @{ Row = 9;  Col = 2; New = "53÷8=" },
    @{ Row = 9;  Col = 3; New = "40÷8=" },
    @{ Row = 9;  Col = 4; New = "63÷5=" },
    @{ Row = 9;  Col = 5; New = "46÷9=" },

    @{ Row = 13; Col = 1; New = "66÷9=" },
    @{ Row = 13; Col = 2; New = "27÷3=" },
    @{ Row = 13; Col = 3; New = "14÷7=" },
    @{ Row = 13; Col = 4; New = "25÷4=" },
    @{ Row = 13; Col = 5; New = "64÷6=" },

    @{ Row = 17; Col = 1; New = "58÷3=" },
    @{ Row = 17; Col = 2; New = "66÷6=" },
    @{ Row = 17; Col = 3; New = "28÷9=" },
    @{ Row = 17; Col = 4; New = "53÷5=" },
    @{ Row = 17; Col = 5; New = "84÷5=" }
)

foreach ($e in $edits) {
    $t.Cell($e.Row, $e.Col).Range.Text = $e.New
}
